# Update the FC Barcelona fixture list:
#  - insert a new match row (Napoli SSC) between the Rayo Vallecano and
#    Athletic Club de Bilbao rows
#  - refresh several of the score/attendance figures in column C
#
# Every value in column C is stored as text in the source workbook (it
# round-trips as a shared string, not a number), so we force column C to
# Text format before writing the new figures and clear the formatting
# again afterwards so no visible style change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "FC Barcelona vs Napoli SSC" above the Athletic
# Club de Bilbao row (current row 3); this pushes the rest of the table
# down by one row.
$ws.Rows(3).Insert()

$colC = $ws.Columns("C")
$colC.NumberFormat = "@"

$ws.Range("A3").Value = "FC Barcelona vs Napoli SSC"
$ws.Range("B3").Value = "February 17th 2022"
$ws.Range("C3").Value = "54"

# FC Barcelona vs Rayo Vallecano
$ws.Range("C2").Value = "75"

# FC Barcelona vs Athletic Club de Bilbao (shifted down to row 4)
$ws.Range("C4").Value = "79"

# FC Barcelona vs Osasuna (shifted down to row 5)
$ws.Range("C5").Value = "79"

# FC Barcelona vs Cadiz CF (shifted down to row 6)
$ws.Range("C6").Value = "75"

# FC Barcelona vs Celta de Vigo (shifted down to row 8)
$ws.Range("C8").Value = "84"

# FC Barcelona vs Villareal CF (shifted down to row 9)
$ws.Range("C9").Value = "81"

$colC.ClearFormats()
